# Regenerate Excel, and add back to list of files we check
#
# 1) Rename sheets (shorten a few list-sheet names, fix the main sheet title)
# 2) Freeze the header row on the main sheet
# 3) Add errorTitle/error text to every data validation rule, and repoint the
#    list-based validations at the renamed lookup sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Sheet renames -------------------------------------------------
$wsMain = $wb.Worksheets.Item("Export this as TSV")
$wsMain.Name = "Export as TSV"
$wb.Worksheets.Item("acquisition_ins-ent_vendor list").Name  = "acquisition_in...nt_vendor list"
$wb.Worksheets.Item("acquisition_ins-ment_model list").Name  = "acquisition_in...ent_model list"
$wb.Worksheets.Item("preparation_ins-ent_vendor list").Name  = "preparation_in...nt_vendor list"
$wb.Worksheets.Item("preparation_ins-ment_model list").Name  = "preparation_in...ent_model list"

# --- 2) Freeze header row on the main sheet ---------------------------
$wsMain.Activate()
$wsMain.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Rewrite data validations with error titles/messages -----------
# NOTE: this interpreter's function dispatch only binds POSITIONAL
# parameters reliably (named args / default values are not bound), so
# every parameter below is passed positionally, in order, every time.
function Set-Validation {
    param($Sqref, $Type, $Formula1, $Formula2, $ErrorTitle, $ErrorMessage)
    $rng = $wsMain.Range($Sqref)
    $rng.Validation.Delete()
    if ($Formula2) {
        $rng.Validation.Add($Type, 1, 1, $Formula1, $Formula2)
    } else {
        $rng.Validation.Add($Type, 1, 1, $Formula1)
    }
    $rng.Validation.ErrorTitle = $ErrorTitle
    $rng.Validation.ErrorMessage = $ErrorMessage
}

Set-Validation 'I2:I1048576' 3 '''assay_category list''!$A$1:$A$1' $null `
    'Value must come from list' 'Value must be one of: imaging.'

Set-Validation 'J2:J1048576' 3 '''assay_type list''!$A$1:$A$1' $null `
    'Value must come from list' 'Value must be one of: CODEX.'

Set-Validation 'K2:K1048576' 3 '''analyte_class list''!$A$1:$A$1' $null `
    'Value must come from list' 'Value must be one of: protein.'

Set-Validation 'L2:L1048576' 3 '"TRUE,FALSE"' $null `
    'Not a boolean' 'The values in this column must be "TRUE" or "FALSE".'

Set-Validation 'M2:M1048576' 3 '''acquisition_in...nt_vendor list''!$A$1:$A$2' $null `
    'Value must come from list' 'Value must be one of: Keyence / Zeiss.'

Set-Validation 'N2:N1048576' 3 '''acquisition_in...ent_model list''!$A$1:$A$3' $null `
    'Value must come from list' 'Value must be one of: BZ-X800 / BZ-X710 / Axio Observer Z1.'

Set-Validation 'O2:O1048576' 2 '-1e+307' '1e+307' `
    'Not a number' 'The values in this column must be numbers.'

Set-Validation 'P2:P1048576' 3 '''resolution_x_unit list''!$A$1:$A$3' $null `
    'Value must come from list' 'Value must be one of: mm / um / nm.'

Set-Validation 'Q2:Q1048576' 2 '-1e+307' '1e+307' `
    'Not a number' 'The values in this column must be numbers.'

Set-Validation 'R2:R1048576' 3 '''resolution_y_unit list''!$A$1:$A$3' $null `
    'Value must come from list' 'Value must be one of: mm / um / nm.'

Set-Validation 'S2:S1048576' 2 '-1e+307' '1e+307' `
    'Not a number' 'The values in this column must be numbers.'

Set-Validation 'T2:T1048576' 3 '''resolution_z_unit list''!$A$1:$A$3' $null `
    'Value must come from list' 'Value must be one of: mm / um / nm.'

Set-Validation 'U2:U1048576' 3 '''preparation_in...nt_vendor list''!$A$1:$A$1' $null `
    'Value must come from list' 'Value must be one of: CODEX.'

Set-Validation 'V2:V1048576' 3 '''preparation_in...ent_model list''!$A$1:$A$2' $null `
    'Value must come from list' 'Value must be one of: version 1 robot / prototype robot - Stanford/Nolan Lab.'

Set-Validation 'W2:W1048576' 1 '-2147483647' '2147483647' `
    'Not an integer' 'The values in this column must be integers.'

Set-Validation 'X2:X1048576' 1 '-2147483647' '2147483647' `
    'Not an integer' 'The values in this column must be integers.'

Set-Validation 'Y2:Y1048576' 1 '-2147483647' '2147483647' `
    'Not an integer' 'The values in this column must be integers.'

$wsMain.Range("A1").Select()
